$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Feeders: drop the "new_department_code" / "new_cost_type_code"
#    columns (C:D) - the feeder rows no longer map to an explicit
#    department/cost type override, only to an event class.
# ------------------------------------------------------------------
$wsFeed = $wb.Worksheets.Item("feeders")
$wsFeed.Activate()
$wsFeed.Range("C1:D1").EntireColumn.Delete()

# ------------------------------------------------------------------
# 2. departments: remove the placeholder "none" / "No associated
#    department" row (row 2) - it is no longer referenced now that
#    feeders don't carry an explicit department override.
# ------------------------------------------------------------------
$wsDept = $wb.Worksheets.Item("departments")
$wsDept.Activate()
$wsDept.Rows(2).Select()
$wsDept.Rows(2).Delete()

# ------------------------------------------------------------------
# 3. cost types: remove the placeholder "none" / "no associated cost
#    type" row (row 2) for the same reason.
# ------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("cost types")
$wsCost.Activate()
$wsCost.Range("J7").Select()
$wsCost.Rows(2).Delete()

# ------------------------------------------------------------------
# 4. Keep the hidden ExternalData_1 / ExternalData_2 query ranges
#    (used by the "dept"/"costtype" Power Query tables) in sync with
#    the now-shorter department/cost-type tables.
# ------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "departments!ExternalData_1") {
        $n.RefersTo = "=departments!`$A`$1:`$B`$61"
    }
    if ($n.Name -eq "cost types!ExternalData_2") {
        $n.RefersTo = "='cost types'!`$A`$1:`$B`$35"
    }
}

# ------------------------------------------------------------------
# Leave the workbook focused back on the feeders sheet, where the
# cursor ends up after the column edit.
# ------------------------------------------------------------------
$wsFeed.Activate()
$wsFeed.Range("I13").Select()
